$wb = $excel.ActiveWorkbook

# The data change is on the "Span 6 RRC 0.6 Tukey" sheet (sheet4):
#   D4: 39 -> 37  (dependent formulas in E4 and E11 recalculate automatically)
$ws = $wb.Worksheets.Item("Span 6 RRC 0.6 Tukey")
$ws.Activate()
$ws.Range("D4").Value = 37

# Update the selection to match the saved view (E6) on that sheet
$ws.Range("E6").Select()
